$wb = $excel.ActiveWorkbook

# ---- Step 1: Insert the new "2022-Q4" sheet before the current "2022-Q3" sheet ----
$beforeSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q4"

# Re-fetch fresh sheet references by name (anchors can go stale after a structural insert)
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Item("2022-Q4")

# ---- Step 2: Seed header row + index column A from the "2022-Q3" template (copies value + style/format) ----
$q3Sheet.Range("A1:H1").Copy($q4Sheet.Range("A1:H1"))
$q3Sheet.Range("A2:A48").Copy($q4Sheet.Range("A2:A48"))

# ---- Step 3: Force text format on the B:G data columns so numeric-looking strings (fund codes, percentages) stay text ----
$q4Sheet.Range("B2:G48").NumberFormat = "@"

# ---- Step 4: Fill in fund data rows 2-48 ----
$q4Sheet.Range("B2").Value = "090018"
$q4Sheet.Range("C2").Value = "大成新锐产业混合"
$q4Sheet.Range("D2").Value = "88.75"
$q4Sheet.Range("E2").Value = "93.33"
$q4Sheet.Range("F2").Value = "7.13"
$q4Sheet.Range("G2").Value = "6.3279"
$q4Sheet.Range("H2").Value = 4
$q4Sheet.Range("B3").Value = "006102"
$q4Sheet.Range("C3").Value = "浙商丰利增强债券"
$q4Sheet.Range("D3").Value = "82.37"
$q4Sheet.Range("E3").Value = "47.70"
$q4Sheet.Range("F3").Value = "4.76"
$q4Sheet.Range("G3").Value = "3.9208"
$q4Sheet.Range("H3").Value = 1
$q4Sheet.Range("B4").Value = "270006"
$q4Sheet.Range("C4").Value = "广发策略优选混合"
$q4Sheet.Range("D4").Value = "48.36"
$q4Sheet.Range("E4").Value = "93.44"
$q4Sheet.Range("F4").Value = "5.86"
$q4Sheet.Range("G4").Value = "2.8339"
$q4Sheet.Range("H4").Value = 4
$q4Sheet.Range("B5").Value = "163415"
$q4Sheet.Range("C5").Value = "兴全商业模式优选混合（LOF）"
$q4Sheet.Range("D5").Value = "108.38"
$q4Sheet.Range("E5").Value = "93.15"
$q4Sheet.Range("F5").Value = "2.33"
$q4Sheet.Range("G5").Value = "2.5253"
$q4Sheet.Range("H5").Value = 10
$q4Sheet.Range("B6").Value = "001300"
$q4Sheet.Range("C6").Value = "大成睿景灵活配置混合A"
$q4Sheet.Range("D6").Value = "32.95"
$q4Sheet.Range("E6").Value = "92.29"
$q4Sheet.Range("F6").Value = "7.13"
$q4Sheet.Range("G6").Value = "2.3493"
$q4Sheet.Range("H6").Value = 4
$q4Sheet.Range("B7").Value = "013435"
$q4Sheet.Range("C7").Value = "大成景气精选六个月持有混合A"
$q4Sheet.Range("D7").Value = "30.45"
$q4Sheet.Range("E7").Value = "91.16"
$q4Sheet.Range("F7").Value = "7.14"
$q4Sheet.Range("G7").Value = "2.1741"
$q4Sheet.Range("H7").Value = 4
$q4Sheet.Range("B8").Value = "013532"
$q4Sheet.Range("C8").Value = "广发安宏回报灵活配置混合E"
$q4Sheet.Range("D8").Value = "24.84"
$q4Sheet.Range("E8").Value = "94.43"
$q4Sheet.Range("F8").Value = "7.69"
$q4Sheet.Range("G8").Value = "1.9102"
$q4Sheet.Range("H8").Value = 2
$q4Sheet.Range("B9").Value = "001761"
$q4Sheet.Range("C9").Value = "广发安宏回报灵活配置混合A"
$q4Sheet.Range("D9").Value = "23.97"
$q4Sheet.Range("E9").Value = "94.43"
$q4Sheet.Range("F9").Value = "7.69"
$q4Sheet.Range("G9").Value = "1.8433"
$q4Sheet.Range("H9").Value = 2
$q4Sheet.Range("B10").Value = "010628"
$q4Sheet.Range("C10").Value = "广发瑞轩三个月定期开放混合"
$q4Sheet.Range("D10").Value = "27.69"
$q4Sheet.Range("E10").Value = "84.10"
$q4Sheet.Range("F10").Value = "5.70"
$q4Sheet.Range("G10").Value = "1.5783"
$q4Sheet.Range("H10").Value = 2
$q4Sheet.Range("B11").Value = "163409"
$q4Sheet.Range("C11").Value = "兴全绿色投资混合（LOF）"
$q4Sheet.Range("D11").Value = "50.87"
$q4Sheet.Range("E11").Value = "91.54"
$q4Sheet.Range("F11").Value = "2.94"
$q4Sheet.Range("G11").Value = "1.4956"
$q4Sheet.Range("H11").Value = 9
$q4Sheet.Range("B12").Value = "001301"
$q4Sheet.Range("C12").Value = "大成睿景灵活配置混合C"
$q4Sheet.Range("D12").Value = "19.30"
$q4Sheet.Range("E12").Value = "92.29"
$q4Sheet.Range("F12").Value = "7.13"
$q4Sheet.Range("G12").Value = "1.3761"
$q4Sheet.Range("H12").Value = 4
$q4Sheet.Range("B13").Value = "002258"
$q4Sheet.Range("C13").Value = "大成国企改革灵活配置混合"
$q4Sheet.Range("D13").Value = "16.71"
$q4Sheet.Range("E13").Value = "93.37"
$q4Sheet.Range("F13").Value = "7.13"
$q4Sheet.Range("G13").Value = "1.1914"
$q4Sheet.Range("H13").Value = 4
$q4Sheet.Range("B14").Value = "014224"
$q4Sheet.Range("C14").Value = "大成聚优成长混合A"
$q4Sheet.Range("D14").Value = "16.55"
$q4Sheet.Range("E14").Value = "90.21"
$q4Sheet.Range("F14").Value = "7.14"
$q4Sheet.Range("G14").Value = "1.1817"
$q4Sheet.Range("H14").Value = 4
$q4Sheet.Range("B15").Value = "470009"
$q4Sheet.Range("C15").Value = "汇添富民营活力混合A"
$q4Sheet.Range("D15").Value = "24.48"
$q4Sheet.Range("E15").Value = "93.48"
$q4Sheet.Range("F15").Value = "4.36"
$q4Sheet.Range("G15").Value = "1.0673"
$q4Sheet.Range("H15").Value = 5
$q4Sheet.Range("B16").Value = "010826"
$q4Sheet.Range("C16").Value = "大成产业趋势混合A"
$q4Sheet.Range("D16").Value = "11.37"
$q4Sheet.Range("E16").Value = "93.99"
$q4Sheet.Range("F16").Value = "7.29"
$q4Sheet.Range("G16").Value = "0.8289"
$q4Sheet.Range("H16").Value = 4
$q4Sheet.Range("B17").Value = "688888"
$q4Sheet.Range("C17").Value = "浙商聚潮产业成长混合A"
$q4Sheet.Range("D17").Value = "9.18"
$q4Sheet.Range("E17").Value = "94.47"
$q4Sheet.Range("F17").Value = "7.62"
$q4Sheet.Range("G17").Value = "0.6995"
$q4Sheet.Range("H17").Value = 3
$q4Sheet.Range("B18").Value = "012519"
$q4Sheet.Range("C18").Value = "大成核心趋势混合A"
$q4Sheet.Range("D18").Value = "9.51"
$q4Sheet.Range("E18").Value = "91.16"
$q4Sheet.Range("F18").Value = "7.14"
$q4Sheet.Range("G18").Value = "0.6790"
$q4Sheet.Range("H18").Value = 4
$q4Sheet.Range("B19").Value = "013531"
$q4Sheet.Range("C19").Value = "浙商聚潮产业成长混合C"
$q4Sheet.Range("D19").Value = "6.24"
$q4Sheet.Range("E19").Value = "94.47"
$q4Sheet.Range("F19").Value = "7.62"
$q4Sheet.Range("G19").Value = "0.4755"
$q4Sheet.Range("H19").Value = 3
$q4Sheet.Range("B20").Value = "010381"
$q4Sheet.Range("C20").Value = "浙商智选价值混合A"
$q4Sheet.Range("D20").Value = "7.03"
$q4Sheet.Range("E20").Value = "91.16"
$q4Sheet.Range("F20").Value = "5.78"
$q4Sheet.Range("G20").Value = "0.4063"
$q4Sheet.Range("H20").Value = 3
$q4Sheet.Range("B21").Value = "010382"
$q4Sheet.Range("C21").Value = "浙商智选价值混合C"
$q4Sheet.Range("D21").Value = "6.53"
$q4Sheet.Range("E21").Value = "91.16"
$q4Sheet.Range("F21").Value = "5.78"
$q4Sheet.Range("G21").Value = "0.3774"
$q4Sheet.Range("H21").Value = 3
$q4Sheet.Range("B22").Value = "013436"
$q4Sheet.Range("C22").Value = "大成景气精选六个月持有混合C"
$q4Sheet.Range("D22").Value = "5.20"
$q4Sheet.Range("E22").Value = "91.16"
$q4Sheet.Range("F22").Value = "7.14"
$q4Sheet.Range("G22").Value = "0.3713"
$q4Sheet.Range("H22").Value = 4
$q4Sheet.Range("B23").Value = "007368"
$q4Sheet.Range("C23").Value = "浙商沪港深精选混合A"
$q4Sheet.Range("D23").Value = "6.05"
$q4Sheet.Range("E23").Value = "92.12"
$q4Sheet.Range("F23").Value = "5.28"
$q4Sheet.Range("G23").Value = "0.3194"
$q4Sheet.Range("H23").Value = 6
$q4Sheet.Range("B24").Value = "002959"
$q4Sheet.Range("C24").Value = "汇添富盈泰灵活配置混合"
$q4Sheet.Range("D24").Value = "5.14"
$q4Sheet.Range("E24").Value = "94.05"
$q4Sheet.Range("F24").Value = "5.03"
$q4Sheet.Range("G24").Value = "0.2585"
$q4Sheet.Range("H24").Value = 9
$q4Sheet.Range("B25").Value = "010827"
$q4Sheet.Range("C25").Value = "大成产业趋势混合C"
$q4Sheet.Range("D25").Value = "3.42"
$q4Sheet.Range("E25").Value = "93.99"
$q4Sheet.Range("F25").Value = "7.29"
$q4Sheet.Range("G25").Value = "0.2493"
$q4Sheet.Range("H25").Value = 4
$q4Sheet.Range("B26").Value = "501065"
$q4Sheet.Range("C26").Value = "汇添富经典成长定期开放混合"
$q4Sheet.Range("D26").Value = "5.18"
$q4Sheet.Range("E26").Value = "89.69"
$q4Sheet.Range("F26").Value = "4.65"
$q4Sheet.Range("G26").Value = "0.2409"
$q4Sheet.Range("H26").Value = 10
$q4Sheet.Range("B27").Value = "014225"
$q4Sheet.Range("C27").Value = "大成聚优成长混合C"
$q4Sheet.Range("D27").Value = "3.26"
$q4Sheet.Range("E27").Value = "90.21"
$q4Sheet.Range("F27").Value = "7.14"
$q4Sheet.Range("G27").Value = "0.2328"
$q4Sheet.Range("H27").Value = 4
$q4Sheet.Range("B28").Value = "012520"
$q4Sheet.Range("C28").Value = "大成核心趋势混合C"
$q4Sheet.Range("D28").Value = "2.45"
$q4Sheet.Range("E28").Value = "91.16"
$q4Sheet.Range("F28").Value = "7.14"
$q4Sheet.Range("G28").Value = "0.1749"
$q4Sheet.Range("H28").Value = 4
$q4Sheet.Range("B29").Value = "009181"
$q4Sheet.Range("C29").Value = "浙商智多兴稳健回报一年持有期混合A"
$q4Sheet.Range("D29").Value = "7.69"
$q4Sheet.Range("E29").Value = "36.88"
$q4Sheet.Range("F29").Value = "1.91"
$q4Sheet.Range("G29").Value = "0.1469"
$q4Sheet.Range("H29").Value = 5
$q4Sheet.Range("B30").Value = "009182"
$q4Sheet.Range("C30").Value = "浙商智多兴稳健回报一年持有期混合C"
$q4Sheet.Range("D30").Value = "4.66"
$q4Sheet.Range("E30").Value = "36.88"
$q4Sheet.Range("F30").Value = "1.91"
$q4Sheet.Range("G30").Value = "0.0890"
$q4Sheet.Range("H30").Value = 5
$q4Sheet.Range("B31").Value = "007369"
$q4Sheet.Range("C31").Value = "浙商沪港深精选混合C"
$q4Sheet.Range("D31").Value = "1.38"
$q4Sheet.Range("E31").Value = "92.12"
$q4Sheet.Range("F31").Value = "5.28"
$q4Sheet.Range("G31").Value = "0.0729"
$q4Sheet.Range("H31").Value = 6
$q4Sheet.Range("B32").Value = "690001"
$q4Sheet.Range("C32").Value = "民生加银品牌蓝筹混合"
$q4Sheet.Range("D32").Value = "1.19"
$q4Sheet.Range("E32").Value = "91.85"
$q4Sheet.Range("F32").Value = "3.43"
$q4Sheet.Range("G32").Value = "0.0408"
$q4Sheet.Range("H32").Value = 6
$q4Sheet.Range("B33").Value = "001762"
$q4Sheet.Range("C33").Value = "广发安宏回报灵活配置混合C"
$q4Sheet.Range("D33").Value = "0.51"
$q4Sheet.Range("E33").Value = "94.43"
$q4Sheet.Range("F33").Value = "7.69"
$q4Sheet.Range("G33").Value = "0.0392"
$q4Sheet.Range("H33").Value = 2
$q4Sheet.Range("B34").Value = "014321"
$q4Sheet.Range("C34").Value = "德邦周期精选混合A"
$q4Sheet.Range("D34").Value = "0.59"
$q4Sheet.Range("E34").Value = "89.92"
$q4Sheet.Range("F34").Value = "5.16"
$q4Sheet.Range("G34").Value = "0.0304"
$q4Sheet.Range("H34").Value = 4
$q4Sheet.Range("B35").Value = "001121"
$q4Sheet.Range("C35").Value = "东方睿鑫热点挖掘灵活配置混合C"
$q4Sheet.Range("D35").Value = "0.66"
$q4Sheet.Range("E35").Value = "90.48"
$q4Sheet.Range("F35").Value = "3.90"
$q4Sheet.Range("G35").Value = "0.0257"
$q4Sheet.Range("H35").Value = 9
$q4Sheet.Range("B36").Value = "001120"
$q4Sheet.Range("C36").Value = "东方睿鑫热点挖掘灵活配置混合A"
$q4Sheet.Range("D36").Value = "0.55"
$q4Sheet.Range("E36").Value = "90.48"
$q4Sheet.Range("F36").Value = "3.90"
$q4Sheet.Range("G36").Value = "0.0214"
$q4Sheet.Range("H36").Value = 9
$q4Sheet.Range("B37").Value = "005351"
$q4Sheet.Range("C37").Value = "汇添富行业整合主题混合A"
$q4Sheet.Range("D37").Value = "0.55"
$q4Sheet.Range("E37").Value = "83.37"
$q4Sheet.Range("F37").Value = "3.84"
$q4Sheet.Range("G37").Value = "0.0211"
$q4Sheet.Range("H37").Value = 6
$q4Sheet.Range("B38").Value = "015191"
$q4Sheet.Range("C38").Value = "汇添富行业整合主题混合D"
$q4Sheet.Range("D38").Value = "0.55"
$q4Sheet.Range("E38").Value = "83.37"
$q4Sheet.Range("F38").Value = "3.84"
$q4Sheet.Range("G38").Value = "0.0211"
$q4Sheet.Range("H38").Value = 6
$q4Sheet.Range("B39").Value = "004244"
$q4Sheet.Range("C39").Value = "东方周期优选灵活配置混合"
$q4Sheet.Range("D39").Value = "0.34"
$q4Sheet.Range("E39").Value = "88.65"
$q4Sheet.Range("F39").Value = "5.79"
$q4Sheet.Range("G39").Value = "0.0197"
$q4Sheet.Range("H39").Value = 6
$q4Sheet.Range("B40").Value = "005161"
$q4Sheet.Range("C40").Value = "华商上游产业股票"
$q4Sheet.Range("D40").Value = "0.55"
$q4Sheet.Range("E40").Value = "88.76"
$q4Sheet.Range("F40").Value = "3.30"
$q4Sheet.Range("G40").Value = "0.0182"
$q4Sheet.Range("H40").Value = 9
$q4Sheet.Range("B41").Value = "011888"
$q4Sheet.Range("C41").Value = "民生加银周期优选混合型证券投资基金A"
$q4Sheet.Range("D41").Value = "0.35"
$q4Sheet.Range("E41").Value = "92.72"
$q4Sheet.Range("F41").Value = "3.46"
$q4Sheet.Range("G41").Value = "0.0121"
$q4Sheet.Range("H41").Value = 8
$q4Sheet.Range("B42").Value = "159990"
$q4Sheet.Range("C42").Value = "银华巨潮小盘价值ETF"
$q4Sheet.Range("D42").Value = "0.78"
$q4Sheet.Range("E42").Value = "97.02"
$q4Sheet.Range("F42").Value = "1.43"
$q4Sheet.Range("G42").Value = "0.0112"
$q4Sheet.Range("H42").Value = 4
$q4Sheet.Range("B43").Value = "011054"
$q4Sheet.Range("C43").Value = "申万菱信安鑫智选混合A"
$q4Sheet.Range("D43").Value = "0.69"
$q4Sheet.Range("E43").Value = "25.20"
$q4Sheet.Range("F43").Value = "1.13"
$q4Sheet.Range("G43").Value = "0.0078"
$q4Sheet.Range("H43").Value = 6
$q4Sheet.Range("B44").Value = "011889"
$q4Sheet.Range("C44").Value = "民生加银周期优选混合型证券投资基金C"
$q4Sheet.Range("D44").Value = "0.09"
$q4Sheet.Range("E44").Value = "92.72"
$q4Sheet.Range("F44").Value = "3.46"
$q4Sheet.Range("G44").Value = "0.0031"
$q4Sheet.Range("H44").Value = 8
$q4Sheet.Range("B45").Value = "014322"
$q4Sheet.Range("C45").Value = "德邦周期精选混合C"
$q4Sheet.Range("D45").Value = "0.05"
$q4Sheet.Range("E45").Value = "89.92"
$q4Sheet.Range("F45").Value = "5.16"
$q4Sheet.Range("G45").Value = "0.0026"
$q4Sheet.Range("H45").Value = 4
$q4Sheet.Range("B46").Value = "960014"
$q4Sheet.Range("C46").Value = "汇添富民营活力混合 O"
$q4Sheet.Range("D46").Value = "0.00"
$q4Sheet.Range("E46").Value = "93.48"
$q4Sheet.Range("F46").Value = "4.36"
$q4Sheet.Range("G46").NumberFormat = "General"
$q4Sheet.Range("G46").Value = 0
$q4Sheet.Range("H46").Value = 5
$q4Sheet.Range("B47").Value = "015190"
$q4Sheet.Range("C47").Value = "汇添富行业整合主题混合C"
$q4Sheet.Range("D47").Value = "0.00"
$q4Sheet.Range("E47").Value = "83.37"
$q4Sheet.Range("F47").Value = "3.84"
$q4Sheet.Range("G47").NumberFormat = "General"
$q4Sheet.Range("G47").Value = 0
$q4Sheet.Range("H47").Value = 6
$q4Sheet.Range("B48").Value = "011055"
$q4Sheet.Range("C48").Value = "申万菱信安鑫智选混合C"
$q4Sheet.Range("D48").Value = "0.00"
$q4Sheet.Range("E48").Value = "25.20"
$q4Sheet.Range("F48").Value = "1.13"
$q4Sheet.Range("G48").NumberFormat = "General"
$q4Sheet.Range("G48").Value = 0
$q4Sheet.Range("H48").Value = 6

# ---- Step 5: Update the "总计" (summary) sheet: shift existing rows down and insert the new 2022-Q4 row at the top ----
$totalSheet = $wb.Worksheets.Item("总计")

# Copy A8 (style template with s=2) down into the brand-new row 9, then overwrite its value
$totalSheet.Range("A8").Copy($totalSheet.Range("A9"))
$totalSheet.Range("A9").Value = 7
$totalSheet.Range("B9").Value = "2021-Q1"
$totalSheet.Range("C9").Value = 5
$totalSheet.Range("D9").Value = 1.46

# Shift rows 8->7..2 down by one (write from bottom up to avoid clobbering source before it is read)
$totalSheet.Range("A8").Value = 6
$totalSheet.Range("B8").Value = "2021-Q2"
$totalSheet.Range("C8").Value = 37
$totalSheet.Range("D8").Value = 14.79
$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2021-Q3"
$totalSheet.Range("C7").Value = 46
$totalSheet.Range("D7").Value = 20.91
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q4"
$totalSheet.Range("C6").Value = 18
$totalSheet.Range("D6").Value = 12.64
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2022-Q1"
$totalSheet.Range("C5").Value = 88
$totalSheet.Range("D5").Value = 46.78
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q2"
$totalSheet.Range("C4").Value = 102
$totalSheet.Range("D4").Value = 59.3
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 58
$totalSheet.Range("D3").Value = 38
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 47
$totalSheet.Range("D2").Value = 37.67

Write-Host "Edit complete"
